$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the sample order-line data (rows 2-4) with generic field-name
# placeholders, matching the header row's semantics, to support splicing
# live database values into the OrderManager template.
$fields = @("sku", "name", "quantity", "cost_per", "total_cost")

for ($row = 2; $row -le 4; $row++) {
    for ($col = 1; $col -le 5; $col++) {
        $ws.Cells.Item($row, $col).Value = $fields[$col - 1]
    }
}
